$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- H column: convert to formulas referencing M column ---
# (values for M must exist logically, but we set M's values below; formulas recalc automatically)
$ws.Range("H1").Formula = "=M1"
$ws.Range("H2:H3").Formula = "=M2"

# --- New columns I:N, filled column by column (top to bottom) to match shared-string order ---
$ws.Range("I1").Value = "Organización_1"
$ws.Range("I2").Value = "Organización_3"
$ws.Range("I3").Value = "Organización_4"

$ws.Range("J1").Value = "Nombre_1 Apellido_1 Segundo_1 | Qa"
$ws.Range("J2").Value = "Nombre_3 Apellido_3 Segundo_3 | Segundo_"
$ws.Range("J3").Value = "Nombre_4 Apellido_4 Segundo_4 | Segun"

$ws.Range("K1").Value = "Demo del acuerdo numero_1"
$ws.Range("K2").Value = "Demo del acuerdo numero_2"
$ws.Range("K3").Value = "Demo del acuerdo numero_3"

$ws.Range("L1").Value = "Iniciado"
$ws.Range("L2").Value = "En Proceso"
$ws.Range("L3").Value = "Concluido"

# M and N filled row by row
$ws.Range("M1").Value = "categoría_1"
$ws.Range("N1").Value = "sigla_1"

$ws.Range("M2").Value = "categoría_2"
$ws.Range("N2").Value = "sigla_2"

$ws.Range("M3").Value = "categoría_3"
$ws.Range("N3").Value = "sigla_3"

# I1 gets the underline style (matches D4's pre-existing style)
$ws.Range("I1").Font.Underline = 2

# --- Empty, styled-only cells (underline style) below the data, through row 10 ---
$ws.Range("M4").Font.Underline = 2
$ws.Range("N4").Font.Underline = 2

$ws.Range("M5").Font.Underline = 2
$ws.Range("N5").Font.Underline = 2

$ws.Range("J6").Font.Underline = 2
$ws.Range("M6").Font.Underline = 2
$ws.Range("N6").Font.Underline = 2

$ws.Range("M7").Font.Underline = 2
$ws.Range("N7").Font.Underline = 2

$ws.Range("M8").Font.Underline = 2
$ws.Range("N8").Font.Underline = 2

$ws.Range("M9").Font.Underline = 2
$ws.Range("N9").Font.Underline = 2

$ws.Range("M10").Font.Underline = 2
$ws.Range("N10").Font.Underline = 2

# --- Column widths for the new columns ---
$ws.Columns.Item(9).ColumnWidth = 12.833333333333332
$ws.Columns.Item(10).ColumnWidth = 37.33333333333333
$ws.Columns.Item(11).ColumnWidth = 24.0

# --- Update selection to match the new active cell ---
$ws.Range("J9").Select() | Out-Null
